$wb = $excel.ActiveWorkbook
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $last)
$ws.Name = "2025-11-10"

$ws.Range("A1").Value = "rank"
$ws.Range("B1").Value = "title"
$ws.Range("C1").Value = "author"
$ws.Range("D1").Value = "latest_episode"
$ws.Range("A1:D1").Font.Bold = $true

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = '新米オッサン冒険者、最強パーティに死ぬほど鍛えられて無敵になる'
$ws.Range("C2").Value = '漫画：荻野ケン 原作：岸馬きらく キャラクター原案：Tea'
$ws.Range("D2").Value = '第72話 後編'
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = '転生コロシアム～最弱スキルで最強の女たちを攻略して奴隷ハーレム作ります～'
$ws.Range("C3").Value = 'zunta(作画) はらわたさいぞう(原作)'
$ws.Range("D3").Value = 'おまけ：ジェットバス'
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = '生徒会にも穴はある！'
$ws.Range("C4").Value = 'むちまろ'
$ws.Range("D4").Value = "第139話`t汐見のしおしお"
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 'とんでもスキルで異世界放浪メシ'
$ws.Range("C5").Value = '赤岸K（漫画） 江口連（原作） 雅（キャラクター原案）'
$ws.Range("D5").Value = '第56話　「大漁まつり」'
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = '時間停止勇者―余命３日の設定じゃ世界を救うには短すぎる―'
$ws.Range("C6").Value = '光永康則'
$ws.Range("D6").Value = '第７１話『扇山停止』⓵'
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = '王子様の友達'
$ws.Range("C7").Value = 'すけろく(著者)'
$ws.Range("D7").Value = '第31話'
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 'いとこのこ'
$ws.Range("C8").Value = 'いぬちく(著者)'
$ws.Range("D8").Value = '第41話'
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = '地元のいじめっ子達に仕返ししようとしたら、別の戦いが始まった。'
$ws.Range("C9").Value = 'マツモトケンゴ'
$ws.Range("D9").Value = '第６８話　ナイトプールの戦いが始まった（１）'
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = '勇者に全部奪われた俺は勇者の母親とパーティを組みました！'
$ws.Range("C10").Value = '久遠まこと(著者) 石のやっさん(原作)'
$ws.Range("D10").Value = '第31話'
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = '元・世界１位のサブキャラ育成日記 ～廃プレイヤー、異世界を攻略中！～'
$ws.Range("C11").Value = '沢村治太郎(原作) 前田理想(漫画) まろ(キャラクター原案)'
$ws.Range("D11").Value = '第80話その1'
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = '実は俺、最強でした？'
$ws.Range("C12").Value = '原作：澄守 彩 漫画：高橋 愛'
$ws.Range("D12").Value = '第130話　ライアスの苦悩・後編'
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = 'このヒーラー、めんどくさい'
$ws.Range("C13").Value = '丹念に発酵(著者)'
$ws.Range("D13").Value = '「コミックス９巻発売記念！　カーラたちが探検中に転移魔法陣を踏んで飛ばされた先を大募集！」結果発表マンガ'
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = '怠惰な悪辱貴族に転生した俺、シナリオをぶっ壊したら規格外の魔力で最凶になった'
$ws.Range("C14").Value = '菊池快晴(原作) 小田童馬(作画) 桑島黎音(キャラクター原案)'
$ws.Range("D14").Value = '第15話'
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = '女友達は頼めば意外とヤらせてくれる'
$ws.Range("C15").Value = 'ろくろ(漫画) 鏡遊(原作)'
$ws.Range("D15").Value = '第26話①'
$ws.Range("A16").Value = 15
$ws.Range("B16").Value = '美人女上司滝沢さん'
$ws.Range("C16").Value = 'やんBARU(著者)'
$ws.Range("D16").Value = '第207話'
$ws.Range("A17").Value = 16
$ws.Range("B17").Value = '異世界魔王と召喚少女の奴隷魔術'
$ws.Range("C17").Value = '原作：むらさきゆきや 漫画：福田直叶 キャラクター原案：鶴崎貴大'
$ws.Range("D17").Value = '第130話　変身してみる（前編）'
$ws.Range("A18").Value = 17
$ws.Range("B18").Value = '勇者パーティを追い出された器用貧乏　～パーティ事情で付与術士をやっていた剣士、万能へと至る～'
$ws.Range("C18").Value = '漫画：よねぞう 原作：都神樹 キャラクター原案：きさらぎゆり'
$ws.Range("D18").Value = '第５４話　勇者を護る器用貧乏（３）'
$ws.Range("A19").Value = 18
$ws.Range("B19").Value = '魔のものたちは企てる'
$ws.Range("C19").Value = '加藤拓弐(原作) ガしガし(作画)'
$ws.Range("D19").Value = '第31話'
$ws.Range("A20").Value = 19
$ws.Range("B20").Value = '不純な彼女達は懺悔しない'
$ws.Range("C20").Value = 'ポロロッカ(著者)'
$ws.Range("D20").Value = '休載イラスト'
$ws.Range("A21").Value = 20
$ws.Range("B21").Value = 'アザミヤコを好きになる'
$ws.Range("C21").Value = 'ユニティコング(原作) ツノニガウ(作画)'
$ws.Range("D21").Value = '第11話前編'
$ws.Range("A22").Value = 21
$ws.Range("B22").Value = 'ライドンキング'
$ws.Range("C22").Value = '馬場康誌'
$ws.Range("D22").Value = '第85話 大統領と宇宙を突く拳'
$ws.Range("A23").Value = 22
$ws.Range("B23").Value = 'バキ外伝 烈海王は異世界転生しても一向にかまわんッッ'
$ws.Range("C23").Value = '板垣恵介 猪原賽 陸井栄史'
$ws.Range("D23").Value = '第83話　ゴーゴン三姉妹'
$ws.Range("A24").Value = 23
$ws.Range("B24").Value = 'リビルドワールド'
$ws.Range("C24").Value = '綾村切人(漫画) ナフセ(原作) 吟(キャラクターデザイン) わいっしゅ(世界観デザイン) cell(メカニックデザイン)'
$ws.Range("D24").Value = '第75話①'
$ws.Range("A25").Value = 24
$ws.Range("B25").Value = '【パクパクですわ】追放されたお嬢様の『モンスターを食べるほど強くなる』スキルは、１食で１レベルアップする前代未聞の最強スキルでした。３日で人類最強になりましたわ～！'
$ws.Range("C25").Value = '島知宏 音速炒飯 有都あらゆる'
$ws.Range("D25").Value = '第２５食　赤スライムのシャーベット、パクパクですわ！（１）'
$ws.Range("A26").Value = 25
$ws.Range("B26").Value = '独身貴族は異世界を謳歌する ～結婚しない男の優雅なおひとりさまライフ～'
$ws.Range("C26").Value = '漫画：駒鳥 ひわ 原作：錬金王 キャラクター原案：三登 いつき'
$ws.Range("D26").Value = '第35話 独身貴族はバーでハイボールを作る（4）'
$ws.Range("A27").Value = 26
$ws.Range("B27").Value = 'よくわからないけれど異世界に転生していたようです'
$ws.Range("C27").Value = '内々けやき あし カオミン'
$ws.Range("D27").Value = '第142話 よくわからないけれど後始末するみたいです（１）'
$ws.Range("A28").Value = 27
$ws.Range("B28").Value = '世界最強の魔女、始めました 〜私だけ『攻略サイト』を見れる世界で自由に生きます〜'
$ws.Range("C28").Value = '戸賀 環 坂木持丸 riritto'
$ws.Range("D28").Value = '第55話②　新居の大掃除をしてみた'
$ws.Range("A29").Value = 28
$ws.Range("B29").Value = '異世界メイドの三ツ星グルメ ～現代ごはん作ったら王宮で大バズリしました～'
$ws.Range("C29").Value = 'モリタ Ｕ４ nima'
$ws.Range("D29").Value = '第14話（２）　春とおぼっちゃまとピクニックランチ（２）'
$ws.Range("A30").Value = 29
$ws.Range("B30").Value = '配信に致命的に向いていない女の子が迷宮で黙々と人助けする配信'
$ws.Range("C30").Value = '下田将也(漫画) 佐藤悪糖(原作) 福きつね(キャラクター原案)'
$ws.Range("D30").Value = '第4話中編'
$ws.Range("A31").Value = 30
$ws.Range("B31").Value = '小林さんちのメイドラゴン'
$ws.Range("C31").Value = 'クール教信者'
$ws.Range("D31").Value = '第153話'
$ws.Range("A32").Value = 31
$ws.Range("B32").Value = '聖者無双'
$ws.Range("C32").Value = '漫画：秋風緋色 原作：ブロッコリーライオン キャラクター原案：sime'
$ws.Range("D32").Value = '第94話　戦乱のドワーフ王国・奴隷の扱い（前半）'
$ws.Range("A33").Value = 32
$ws.Range("B33").Value = '姫様“拷問”の時間です'
$ws.Range("C33").Value = '原作:春原ロビンソン　漫画:ひらけい'
$ws.Range("D33").Value = '拷問156'
$ws.Range("A34").Value = 33
$ws.Range("B34").Value = '田舎で恋は難しい!?'
$ws.Range("C34").Value = 'ねこうめ(著者)'
$ws.Range("D34").Value = '第1話'
$ws.Range("A35").Value = 34
$ws.Range("B35").Value = '治癒魔法の間違った使い方 ~戦場を駆ける回復要員~'
$ws.Range("C35").Value = '九我山レキ(漫画) くろかた(原作) ＫｅＧ(キャラクター原案)'
$ws.Range("D35").Value = '第83話(後編)その2'
$ws.Range("A36").Value = 35
$ws.Range("B36").Value = 'ハズレ枠の【状態異常スキル】で最強になった俺がすべてを蹂躙するまで'
$ws.Range("C36").Value = '鵜吉しょう（作画） 内々けやき（構成） 篠崎 芳（原作） KWKM（キャラクター原案）'
$ws.Range("D36").Value = '第59話　別れと、出立（前編）'
$ws.Range("A37").Value = 36
$ws.Range("B37").Value = '経験値貯蓄でのんびり傷心旅行 ～勇者と恋人に追放された戦士の無自覚ざまぁ～'
$ws.Range("C37").Value = '奏ヨシキ(著者) 徳川レモン(原作) riritto(キャラクターデザイン)'
$ws.Range("D37").Value = '第40話-2'
$ws.Range("A38").Value = 37
$ws.Range("B38").Value = 'バキ外伝　ガイアとシコルスキー　～ときどきノムラ 二人だけど三人暮らし～'
$ws.Range("C38").Value = '板垣恵介 林たかあき'
$ws.Range("D38").Value = '第57話 銃撃戦'
$ws.Range("A39").Value = 38
$ws.Range("B39").Value = 'くらいあの子としたいこと'
$ws.Range("C39").Value = '碇マナツ(著者)'
$ws.Range("D39").Value = '特別編㉑'
$ws.Range("A40").Value = 39
$ws.Range("B40").Value = '理想のヒモ生活'
$ws.Range("C40").Value = '日月ネコ(漫画) 渡辺恒彦（ヒーロー文庫／イマジカインフォス）(原作) 文倉十(キャラクター原案)'
$ws.Range("D40").Value = '第89話　その1'
$ws.Range("A41").Value = 40
$ws.Range("B41").Value = '十年目、帰還を諦めた転移者はいまさら主人公になる'
$ws.Range("C41").Value = '原作：氷純（「十年目、帰還を諦めた転移者はいまさら主人公になる」MFブックス刊） 漫画：しゅーかま キャラクター原案：あんべよしろう'
$ws.Range("D41").Value = '第２０話③'
$ws.Range("A42").Value = 41
$ws.Range("B42").Value = '追放されたチート付与魔術師は 気ままなセカンドライフを謳歌する。'
$ws.Range("C42").Value = '六志麻あさ 業務用餅 kisui'
$ws.Range("D42").Value = '第７４話ー②'
$ws.Range("A43").Value = 42
$ws.Range("B43").Value = '地味子な三葉さんが僕を誘惑する'
$ws.Range("C43").Value = 'はぶらえる(著者)'
$ws.Range("D43").Value = '第12話前半'
$ws.Range("A44").Value = 43
$ws.Range("B44").Value = '落ちこぼれだった兄が実は最強 ～史上最強の勇者は転生し、学園で無自覚に無双する～'
$ws.Range("C44").Value = '村上よしゆき 茨木野 あるてら'
$ws.Range("D44").Value = '【描き下ろしイラスト】休載です！（その4）'
$ws.Range("A45").Value = 44
$ws.Range("B45").Value = '男嫌いな美人姉妹を名前も告げずに助けたら一体どうなる?'
$ws.Range("C45").Value = 'みょん(原作) 司馬淳子(漫画) ぎうにう(キャラクターデザイン)'
$ws.Range("D45").Value = 'コミックス第4巻発売告知'
$ws.Range("A46").Value = 45
$ws.Range("B46").Value = '黄金の経験値'
$ws.Range("C46").Value = '原純(原作) 霜月汐(作画) fixro2n(キャラクター原案)'
$ws.Range("D46").Value = '第19話（後編）'
$ws.Range("A47").Value = 46
$ws.Range("B47").Value = '10年ぶりに再会したクソガキは清純美少女JKに成長していた'
$ws.Range("C47").Value = '緑青黒羽（漫画） 館西夕木（原作） ひげ猫（キャラクター原案）'
$ws.Range("D47").Value = '第7話　ショッピングデート（前編）'
$ws.Range("A48").Value = 47
$ws.Range("B48").Value = '婚約者に裏切られた錬金術師は、独立して『ざまぁ』します　コミック版'
$ws.Range("C48").Value = '漫画/すたひろ 原作/Y.A'
$ws.Range("D48").Value = 'chapter73【38話②】'
$ws.Range("A49").Value = 48
$ws.Range("B49").Value = '俺以外誰も採取できない素材なのに「素材採取率が低い」とパワハラする幼馴染錬金術師と絶縁した専属魔導士、辺境の町でスローライフを送りたい。'
$ws.Range("C49").Value = '狐御前(原作) 西岡知三(作画) ＮＯＣＯ(キャラクター原案)'
$ws.Range("D49").Value = '第27話-1'
$ws.Range("A50").Value = 49
$ws.Range("B50").Value = '塔の管理をしてみよう'
$ws.Range("C50").Value = '盧恩＆雪笠(Friendly Land)(著者) 早秋(原作) 雨神(キャラクター原案)'
$ws.Range("D50").Value = '第94話後編'
$ws.Range("A51").Value = 50
$ws.Range("B51").Value = '宇崎ちゃんは遊びたい！'
$ws.Range("C51").Value = '丈(著者)'
$ws.Range("D51").Value = '第129話'
